$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a Text number format on price cells whose new value would
# otherwise be auto-parsed by Excel as a numeric literal (these are
# plain-text price strings in the source data, e.g. "42.84").
$textForceCells = @(
    'D5',
    'D6',
    'D8',
    'D9',
    'D10',
    'D11',
    'D12',
    'D13',
    'D14',
    'D19',
    'D21',
    'D22',
    'D27',
    'D28',
    'D29',
    'D30',
    'D31',
    'D32',
    'D33',
    'D34',
    'D36',
    'D39',
    'D41',
    'D42',
    'D44',
    'D47',
    'D50',
    'D51'
)
foreach ($cellAddr in $textForceCells) {
    $ws.Range($cellAddr).NumberFormat = '@'
}

# Apply the updated cell values from the crypto-price refresh.
$ws.Range('D2').Value = '35.359.42'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').Value = '1.893.25'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '246.27'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').Value = '0.691'
$ws.Range('E6').Value = '  +1.76%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '42.84'
$ws.Range('E8').Value = '  -0.38%  '
$ws.Range('D9').Value = '0.357'
$ws.Range('E9').Value = '  +3.66%  '
$ws.Range('D10').Value = '56.34'
$ws.Range('E10').Value = '  +9.54%  '
$ws.Range('D11').Value = '0.0758'
$ws.Range('E11').Value = '  +3.46%  '
$ws.Range('D12').Value = '0.0981'
$ws.Range('E12').Value = '  +1.33%  '
$ws.Range('D13').Value = '13.96'
$ws.Range('E13').Value = '  +7.51%  '
$ws.Range('D14').Value = '0.792'
$ws.Range('E14').Value = '  +10.89%  '
$ws.Range('D15').Value = '2.167.70'
$ws.Range('E15').Value = '  +1.10%  '
$ws.Range('E16').Value = '  +3.15%  '
$ws.Range('D17').Value = '1.899.74'
$ws.Range('E17').Value = '  +0.92%  '
$ws.Range('D18').Value = '35.338.63'
$ws.Range('E18').Value = '  +1.39%  '
$ws.Range('D19').Value = '73.58'
$ws.Range('E19').Value = '  +1.16%  '
$ws.Range('D20').Value = '0.0₃0831'
$ws.Range('E20').Value = '  +2.57%  '
$ws.Range('D21').Value = '244.23'
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').Value = '12.95'
$ws.Range('E22').Value = '  +2.38%  '
$ws.Range('E23').Value = '  +6.24%  '
$ws.Range('E24').Value = '  +7.31%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('D27').Value = '167.17'
$ws.Range('E27').Value = '  +2.12%  '
$ws.Range('D28').Value = '8.59'
$ws.Range('E28').Value = '  +2.14%  '
$ws.Range('D29').Value = '18.31'
$ws.Range('E29').Value = '  +1.25%  '
$ws.Range('D30').Value = '0.128'
$ws.Range('E30').Value = '  +1.69%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '4.36'
$ws.Range('E31').Value = '  +3.79%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.0603'
$ws.Range('E32').Value = '  +5.41%  '
$ws.Range('D33').Value = '4.24'
$ws.Range('E33').Value = '  +1.82%  '
$ws.Range('D34').Value = '1.87'
$ws.Range('E34').Value = '  +26.13%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').Value = '1.48'
$ws.Range('E36').Value = '  -16.39%  '
$ws.Range('E37').Value = '  +2.40%  '
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('D39').Value = '0.0732'
$ws.Range('E39').Value = '  +9.96%  '
$ws.Range('E40').Value = '  +6.96%  '
$ws.Range('D41').Value = '98.95'
$ws.Range('E41').Value = '  +1.36%  '
$ws.Range('D42').Value = '16.88'
$ws.Range('E42').Value = '  -1.39%  '
$ws.Range('E43').Value = '  +0.42%  '
$ws.Range('D44').Value = '13.77'
$ws.Range('E44').Value = '  +16.69%  '
$ws.Range('D45').Value = '1.327.93'
$ws.Range('E45').Value = '  +3.33%  '
$ws.Range('E46').Value = '  +1.80%  '
$ws.Range('D47').Value = '0.0811'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('E48').Value = '  +0.87%  '
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('D50').Value = '6.38'
$ws.Range('E50').Value = '  +0.87%  '
$ws.Range('D51').Value = '42.53'
$ws.Range('E51').Value = '  -0.17%  '
